$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1190320826869504
    "C2" = 0.04071648406533734
    "D2" = 0.7527432677738641
    "E2" = 0.4942365360607697
    "G2" = 1.406728370586922

    "B3" = 0.2917716402565462
    "C3" = 0.04071648406533734
    "D3" = 0.1494219747398047
    "E3" = 0.4942365360607697
    "G3" = 0.9761466351224579

    "B4" = 3.286832544864788
    "C4" = 1.655778082260271
    "D4" = 0.1494219747398047
    "E4" = 10.19245300693656
    "G4" = 15.28448560880142

    "B5" = 3.286832544864788
    "C5" = 1.655778082260271
    "D5" = 0.7527432677738641
    "E5" = 0.4942365360607697
    "G5" = 6.189590430959694

    "B6" = 0.6606524410359556
    "C6" = 0.002571899574220771
    "D6" = 0.1494219747398047
    "E6" = 0.4942365360607697
    "G6" = 1.306882851410751

    "B7" = 3.286832544864788
    "C7" = 1.655778082260271
    "D7" = 0.7527432677738641
    "E7" = 0.4942365360607697
    "G7" = 6.189590430959694

    "B8" = 1.455362044514542
    "C8" = 0.306821227259698
    "D8" = 0.7527432677738641
    "E8" = 10.19245300693656
    "G8" = 12.70737954648466
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
